# Switch the presentation's applied theme color scheme from the
# "Integral" (Red Violet) design to the default "Office Theme" palette.
#
# The deck ships two theme parts: theme1.xml (the live theme, used by
# the slide master / all layouts / all slides - currently "Integral",
# "Red Violet" colours) and theme2.xml (an unused spare theme, the
# default "Office Theme" palette, only linked from the notes master).
# The edit swaps which palette is which: the live theme becomes the
# plain "Office Theme" colours.
#
# Colours are pushed through the live ThemeColorScheme on the slide
# master (the only theme surface PowerPoint's automation model exposes
# for writes) using the standard VBA-style RGB() packing
# (R + G*256 + B*65536).

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

$tcs.Colors(1).RGB  = 0           # dk1      000000
$tcs.Colors(2).RGB  = 16777215    # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388     # dk2      44546A
$tcs.Colors(4).RGB  = 15132391    # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939    # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501     # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845    # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407       # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308    # accent5  4472C4
$tcs.Colors(10).RGB = 4697456     # accent6  70AD47
$tcs.Colors(11).RGB = 12673797    # hlink    0563C1
$tcs.Colors(12).RGB = 7491477     # folHlink 954F72
